# Generate Report for Handoff
# Adds a new handed-off file ("05297bcb-...md") as a new row (row 3) to the
# Overview / zh-cn / de-de worksheets, mirroring the existing row 2 layout.

$wb = $excel.ActiveWorkbook

$newMdName   = "05297bcb-658c-4299-9942-7266161d82d9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdDisplay = "e2e\" + $newMdName
$newMdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e28d653a21273af205cb5bb26327ef5015752208/e2e/" + $newMdName

$newZhXlf = "05297bcb-658c-4299-9942-7266161d82d9oooooooooooooooooooooooooooooooooooooooo.b1694ff1d25f204fda7d8a302255c5536ae7b9e3.zh-cn.xlf"
$newDeXlf = "05297bcb-658c-4299-9942-7266161d82d9oooooooooooooooooooooooooooooooooooooooo.b1694ff1d25f204fda7d8a302255c5536ae7b9e3.de-de.xlf"

$statusText = "Ready for handoff"
$handoffDate = "2016-08-15 10:29:30"
$xliffDate   = "2016-08-15 10:29:25"

# ---------------------------------------------------------------------------
# Sheet "Overview": append row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", $newMdDisplay)
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G3").Value = "'" + $handoffDate

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn": append row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newMdUrl, "", "", $newMdName)
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = "'" + $handoffDate
$wsZh.Range("I3").Value = "'"
$wsZh.Range("J3").Value = "'"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "'0001-01-01 00:00:00"
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = "'"

$wsZh.Columns.Item(3).ColumnWidth = 16.33

$tblZh = $wsZh.ListObjects.Item(1)
$tblZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de": append row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newMdUrl, "", "", $newMdName)
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = "'" + $handoffDate
$wsDe.Range("I3").Value = "'"
$wsDe.Range("J3").Value = "'"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "'0001-01-01 00:00:00"
$wsDe.Range("L3").Value = "'"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = "'"
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = "'"

$wsDe.Columns.Item(3).ColumnWidth = 16.33

$tblDe = $wsDe.ListObjects.Item(1)
$tblDe.Resize($wsDe.Range("A1:P3"))
